$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 7: "TextBox 47" — "[addressbook is modified]" -> "[command modifies addressbook]"
$sh7 = $s.Shapes.Item(7)
$tr7 = $sh7.TextFrame.TextRange
# Edit from the end backwards so earlier offsets remain valid
$tr7.Characters(13, 13).Text = "]"
$tr7.Characters(1, 1).Text = "[command modifies "

# Resize/move TextBox 47 per the updated layout
# (point values tuned so the Single-precision COM round-trip lands on the
# exact target EMU offsets: x=4828265 y=1865986 cx=1472017 cy=923714)
$sh7.Left = 380.1783661417323
$sh7.Top = 146.92805118110238
$sh7.Width = 115.9068503937008
$sh7.Height = 72.73340551181103

# --- Shape 8: "Rectangle: Rounded Corners 50" — "Add ..." -> "Save ..." and split trailing run
$sh8 = $s.Shapes.Item(8)
$tr8 = $sh8.TextFrame.TextRange
# Split the trailing " and clear redundant states" run into " " + "and clear redundant states"
# (done first, while character offsets still match the original "Add ..." text)
$tr8.Characters(41, 26).Font.Size = 18.01
$tr8.Characters(1, 4).Text = "Save "
